$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.986.67"
$ws.Range("E2").Value = "  -3.03%  "
$ws.Range("D3").Value = "2.825.90"
$ws.Range("E3").Value = "  -3.61%  "
$ws.Range("E4").Value = "  -0.24%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.26"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -5.52%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.66"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -7.60%  "
$ws.Range("E7").Value = "  -0.12%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -5.92%  "
$ws.Range("D9").Value = "2.823.88"
$ws.Range("E9").Value = "  -3.66%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -7.52%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.90"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("E12").Value = "  -3.85%  "
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "3.314.61"
$ws.Range("E14").Value = "  -4.06%  "
$ws.Range("D15").Value = "59.156.46"
$ws.Range("E15").Value = "  -3.12%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.39"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -7.83%  "
$ws.Range("D17").Value = "2.838.43"
$ws.Range("E17").Value = "  -3.50%  "
$ws.Range("E18").Value = "  -6.39%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.69"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -7.08%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.98"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -6.83%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "346.80"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -5.81%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.23"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -5.46%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  -0.72%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.88"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -3.03%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.425"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -7.34%  "
$ws.Range("E27").Value = "  -7.23%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +0.34%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -6.84%  "
$ws.Range("D30").Value = "0.0₃0795"
$ws.Range("E30").Value = "  -10.29%  "
$ws.Range("E31").Value = "  -0.06%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.60"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -5.31%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.96"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -5.42%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.37"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -5.13%  "
$ws.Range("E35").Value = "  -7.40%  "
$ws.Range("E36").Value = "  -7.42%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.899"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -12.18%  "
$ws.Range("E38").Value = "  -9.71%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.70"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -3.71%  "
$ws.Range("D40").Value = "2.213.01"
$ws.Range("E40").Value = "  -6.83%  "
$ws.Range("E41").Value = "  -3.97%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.51"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -7.53%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  -4.21%  "
$ws.Range("E45").Value = "  -10.71%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.07"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -10.95%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.35"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  -5.38%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0886"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -5.48%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.53"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -10.54%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.28"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -8.49%  "
